$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1264.7262
$ws.Range("I15").Value = 1264.7262
$ws.Range("K15").Value = 3794.1786
$ws.Range("M15").Value = -3625.1786
$ws.Range("H47").Value = 50000
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H112").Value = 26317180
$ws.Range("J112").Value = 1498.0294
$ws.Range("L112").Value = 4494.0882
$ws.Range("N112").Value = -6710.0882
$ws.Range("H129").Value = 822.12
$ws.Range("J129").Value = 899.2558
$ws.Range("L129").Value = 2697.7674
$ws.Range("N129").Value = -12697.7674
$ws.Range("H137").Value = 1324786.9
$ws.Range("I137").Value = 2071784.6
$ws.Range("K137").Value = 6215353.800000001
$ws.Range("M137").Value = -6212803.800000001
$ws.Range("H138").Value = 7628.53
$ws.Range("I138").Value = 2317.5
$ws.Range("J138").Value = 7849.8228
$ws.Range("K138").Value = 6952.5
$ws.Range("L138").Value = 23549.4684
$ws.Range("M138").Value = -1812.5
$ws.Range("N138").Value = -33829.4684
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1154.8
$ws.Range("I2").Value = 1154.8
$ws.Range("K2").Value = 1154.8
$ws.Range("M2").Value = -1041.8
$ws.Range("H32").Value = 4634.9077
$ws.Range("I32").Value = 3908.3208
$ws.Range("J32").Value = 7844
$ws.Range("K32").Value = 3908.3208
$ws.Range("L32").Value = 7844
$ws.Range("M32").Value = -3621.3208
$ws.Range("N32").Value = -8418
$ws.Range("H61").Value = 4752.6895
$ws.Range("I61").Value = 1276.619
$ws.Range("J61").Value = 13877.375
$ws.Range("K61").Value = 1276.619
$ws.Range("L61").Value = 13877.375
$ws.Range("M61").Value = -1064.619
$ws.Range("N61").Value = -14301.375
$ws.Range("H74").Value = 4563.517
$ws.Range("I74").Value = 5102.6313
$ws.Range("J74").Value = 3539.2
$ws.Range("K74").Value = 5102.6313
$ws.Range("L74").Value = 3539.2
$ws.Range("M74").Value = -4228.6313
$ws.Range("N74").Value = -5287.2
$ws.Range("H77").Value = 4563.517
$ws.Range("I77").Value = 5102.6313
$ws.Range("J77").Value = 3539.2
$ws.Range("K77").Value = 25513.1565
$ws.Range("L77").Value = 17696
$ws.Range("M77").Value = -21145.1565
$ws.Range("N77").Value = -26432
$ws.Range("H116").Value = 1154.8
$ws.Range("I116").Value = 1154.8
$ws.Range("K116").Value = 1154.8
$ws.Range("M116").Value = 1139.2
$ws.Range("H122").Value = 3559.7778
$ws.Range("I122").Value = 1477.5
$ws.Range("K122").Value = 4432.5
$ws.Range("M122").Value = -1982.5
$ws.Range("H132").Value = 1762.4237
$ws.Range("I132").Value = 1167.5526
$ws.Range("J132").Value = 2838.8572
$ws.Range("K132").Value = 3502.6578
$ws.Range("L132").Value = 8516.571599999999
$ws.Range("M132").Value = -972.6578
$ws.Range("N132").Value = -13576.5716
$ws.Range("H136").Value = 4752.6895
$ws.Range("I136").Value = 1276.619
$ws.Range("J136").Value = 13877.375
$ws.Range("K136").Value = 3829.857
$ws.Range("L136").Value = 41632.125
$ws.Range("M136").Value = -1279.857
$ws.Range("N136").Value = -46732.125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1154.8
$ws.Range("I3").Value = 1154.8
$ws.Range("K3").Value = 1154.8
$ws.Range("M3").Value = -1040.8
$ws.Range("H134").Value = 2713.065
$ws.Range("I134").Value = 1031.5084
$ws.Range("J134").Value = 8224.833000000001
$ws.Range("K134").Value = 3094.5252
$ws.Range("L134").Value = 24674.499
$ws.Range("M134").Value = -559.5252
$ws.Range("N134").Value = -29744.499
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3233.84
$ws.Range("I31").Value = 1230.4706
$ws.Range("K31").Value = 1230.4706
$ws.Range("M31").Value = -935.4706000000001
$ws.Range("H34").Value = 3233.84
$ws.Range("I34").Value = 1230.4706
$ws.Range("K34").Value = 1230.4706
$ws.Range("M34").Value = -1028.4706
$ws.Range("H58").Value = 2594.6619
$ws.Range("I58").Value = 1715.3103
$ws.Range("J58").Value = 7694.9
$ws.Range("K58").Value = 1715.3103
$ws.Range("L58").Value = 7694.9
$ws.Range("M58").Value = -1512.3103
$ws.Range("N58").Value = -8100.9
$ws.Range("H132").Value = 3359.3333
$ws.Range("I132").Value = 2840
$ws.Range("J132").Value = 7124.5
$ws.Range("K132").Value = 8520
$ws.Range("L132").Value = 21373.5
$ws.Range("M132").Value = -5990
$ws.Range("N132").Value = -26433.5
$ws.Range("H134").Value = 1573.4595
$ws.Range("I134").Value = 965.4761999999999
$ws.Range("J134").Value = 2371.4375
$ws.Range("K134").Value = 2896.4286
$ws.Range("L134").Value = 7114.3125
$ws.Range("M134").Value = -361.4285999999997
$ws.Range("N134").Value = -12184.3125
$ws.Range("H136").Value = 2594.6619
$ws.Range("I136").Value = 1715.3103
$ws.Range("J136").Value = 7694.9
$ws.Range("K136").Value = 5145.9309
$ws.Range("L136").Value = 23084.7
$ws.Range("M136").Value = -2595.9309
$ws.Range("N136").Value = -28184.7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1788.421
$ws.Range("I3").Value = 1434.1177
$ws.Range("K3").Value = 4302.3531
$ws.Range("M3").Value = -4190.3531
$ws.Range("H60").Value = 12849.889
$ws.Range("J60").Value = 16485.143
$ws.Range("L60").Value = 49455.429
$ws.Range("N60").Value = -49957.429
$ws.Range("H131").Value = 6098357
$ws.Range("I131").Value = 62500416
$ws.Range("J131").Value = 836.7027
$ws.Range("K131").Value = 187501248
$ws.Range("L131").Value = 2510.1081
$ws.Range("M131").Value = -187496208
$ws.Range("N131").Value = -12590.1081
$ws.Range("H133").Value = 3412.9167
$ws.Range("I133").Value = 3036.125
$ws.Range("K133").Value = 9108.375
$ws.Range("M133").Value = -4048.375
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3042.5
$ws.Range("I122").Value = 2358.6875
$ws.Range("J122").Value = 3824
$ws.Range("K122").Value = 7076.0625
$ws.Range("L122").Value = 11472
$ws.Range("M122").Value = -4626.0625
$ws.Range("N122").Value = -16372
$ws.Range("H132").Value = 2277.138
$ws.Range("I132").Value = 1121.5883
$ws.Range("J132").Value = 3914.1667
$ws.Range("K132").Value = 3364.7649
$ws.Range("L132").Value = 11742.5001
$ws.Range("M132").Value = -834.7648999999997
$ws.Range("N132").Value = -16802.5001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4605.6665
$ws.Range("I122").Value = 1908.5
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 5725.5
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -3275.5
$ws.Range("N122").Value = -34900
$ws.Range("H136").Value = 2947.318
$ws.Range("J136").Value = 5033.0557
$ws.Range("L136").Value = 15099.1671
$ws.Range("N136").Value = -20199.1671
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1819.0476
$ws.Range("I132").Value = 1263.4482
$ws.Range("J132").Value = 3058.4614
$ws.Range("K132").Value = 3790.3446
$ws.Range("L132").Value = 9175.3842
$ws.Range("M132").Value = -1260.3446
$ws.Range("N132").Value = -14235.3842
$ws.Range("H136").Value = 2982.4055
$ws.Range("I136").Value = 2086.8462
$ws.Range("J136").Value = 5099.1816
$ws.Range("K136").Value = 6260.5386
$ws.Range("L136").Value = 15297.5448
$ws.Range("M136").Value = -3710.5386
$ws.Range("N136").Value = -20397.5448
